# Add data for 2021-12-31: the "through December 22" running tally becomes
# "through December 23" — rename the sheet, update the column header text,
# and bump the affected neighborhood/month cell counts (+ a few brand-new
# non-zero cells) for the December-2021 column plus a handful of scattered
# updates elsewhere in the grid.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sheet name & running-total header text ---
$ws.Name = "Through 2021-12-23"
$ws.Range("B1").Value = "December 2021 (through December 23)"

# --- Modified existing cell values ---
$ws.Range("B3").Value   = 9
$ws.Range("B4").Value   = 6
$ws.Range("B6").Value   = 4
$ws.Range("B7").Value   = 7
$ws.Range("N7").Value   = 8
$ws.Range("AL7").Value  = 12
$ws.Range("N9").Value   = 4
$ws.Range("N13").Value  = 3
$ws.Range("N15").Value  = 7
$ws.Range("BJ21").Value = 2
$ws.Range("AX22").Value = 4
$ws.Range("BJ24").Value = 5
$ws.Range("BJ25").Value = 2
$ws.Range("B49").Value  = 3
$ws.Range("B61").Value  = 3
$ws.Range("N99").Value  = 2

# --- Newly populated (previously empty) cells ---
$ws.Range("BJ3").Value  = 1
$ws.Range("Z12").Value  = 1
$ws.Range("BV16").Value = 1
$ws.Range("Z18").Value  = 1
$ws.Range("AX29").Value = 1
$ws.Range("AX32").Value = 1
$ws.Range("AL37").Value = 1
$ws.Range("AL70").Value = 1
$ws.Range("Z88").Value  = 1
$ws.Range("BV89").Value = 1
$ws.Range("N93").Value  = 2

Write-Output "applied carjacking-by-neighborhood-by-month update for 2021-12-23"
